# "Generate Report for Archive"
#
# The localization status report is being regenerated: the status text for
# the fc632082...md file flips from "Ready for handoff" to "In Translation"
# on every sheet that surfaces it (Overview!E2/F2, zh-cn!C2, de-de!C2 all
# share the same status string), and the "Status" column is re-sized to fit
# the new (shorter) text on the Overview sheet (columns E & F) and on the
# per-locale sheets (column C on zh-cn and de-de).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Update the status text everywhere it appears, across every worksheet.
# (Whole-cell match so we only touch cells that are exactly the old status,
# not merely contain it as a substring.)
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace($oldStatus, $newStatus, [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
}

# Re-fit the "Status" columns now that the text is shorter.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5   # column E (zh-cn status)
$wsOverview.Columns.Item(6).ColumnWidth = 12.5   # column F (de-de status)

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
